# notation_translation_table.xlsx - aligning notation and tidying code
# - add gamma ("γ") to the topicmodels notation column (E) next to the
#   existing Blei/LDATS gamma row, and to a new row beneath it
# - add the new LDATS "standardized concentration parameter" row (ybar /
#   upsilon-with-macron) describing the organism-level community identity
#   probability parameter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gamma       = [char]0x03B3   # γ
$upsilon     = [char]0x03C5   # υ
$upsilonBar  = [char]0x1FE1   # ῡ

# --- row 22: align the "topicmodels" notation (E) with the existing Blei
#     notation (C22 is already "γ") -------------------------------------
$ws.Range("E22").Value = $gamma
$ws.Range("E22").Font.Name = "Calibri"

# --- row 24 (new): same alignment cell in column E, plus the new LDATS
#     notation/meaning pair in columns G/H --------------------------------
$ws.Range("E24").Value = $gamma
$ws.Range("E24").Font.Name = "Calibri"

$ws.Range("G24").Value = $upsilonBar
$ws.Range("G24").Font.Name = "Calibri"

$meaning = "standardized concentration parmater " + $upsilon + ", i.e. organism-level community identity probability"
$ws.Range("H24").Value = $meaning
# format the Greek letter and everything after it (", i.e. ...") to match
# the authored file's distinct (non-theme) Calibri run
$ws.Range("H24").Characters(37, 53).Font.Name = "Calibri"

# --- selection moves to G26 after the edits, as in the authored file ----
$null = $ws.Range("G26").Select()
